$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 0.747119
$ws.Range("H2").Value = 2.241357
$ws.Range("I2").Value = 0.03096954854571248
$ws.Range("J2").Value = 0.03096954854571248
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 0.7065936666666666
$ws.Range("N2").Value = 2.119781
$ws.Range("O2").Value = 0.005187843618793344
$ws.Range("P2").Value = 0.005187843618793344
$ws.Range("Q2").Value = 0.5279095536463333
$ws.Range("R2").Value = 4.751185982816999
$ws.Range("S2").Value = 0.0001606651747997852
$ws.Range("T2").Value = 0.0001606651747997852

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 0.747119
$ws.Range("H3").Value = 2.241357
$ws.Range("I3").Value = 0.03096954854571248
$ws.Range("J3").Value = 0.03096954854571248
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 111.9320066666667
$ws.Range("N3").Value = 335.79602
$ws.Range("O3").Value = 0.8218100075305903
$ws.Range("P3").Value = 0.8218100075305903
$ws.Range("Q3").Value = 83.62652888879333
$ws.Range("R3").Value = 752.6387599991399
$ws.Range("S3").Value = 0.02545108492357095
$ws.Range("T3").Value = 0.02545108492357095

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 0.747119
$ws.Range("H4").Value = 2.241357
$ws.Range("I4").Value = 0.03096954854571248
$ws.Range("J4").Value = 0.03096954854571248
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 23.563205
$ws.Range("N4").Value = 70.689615
$ws.Range("O4").Value = 0.1730021488506163
$ws.Range("P4").Value = 0.1730021488506163
$ws.Range("Q4").Value = 17.604518156395
$ws.Range("R4").Value = 158.440663407555
$ws.Range("S4").Value = 0.005357798447341738
$ws.Range("T4").Value = 0.005357798447341739

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 19.74619233333334
$ws.Range("H5").Value = 59.23857700000001
$ws.Range("I5").Value = 0.8185184181638298
$ws.Range("J5").Value = 0.8185184181638298
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 0.7065936666666666
$ws.Range("N5").Value = 2.119781
$ws.Range("O5").Value = 0.005187843618793344
$ws.Range("P5").Value = 0.005187843618793344
$ws.Range("Q5").Value = 13.95253444351522
$ws.Range("R5").Value = 125.572809991637
$ws.Range("S5").Value = 0.004246345552536046
$ws.Range("T5").Value = 0.004246345552536046

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 19.74619233333334
$ws.Range("H6").Value = 59.23857700000001
$ws.Range("I6").Value = 0.8185184181638298
$ws.Range("J6").Value = 0.8185184181638298
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 111.9320066666667
$ws.Range("N6").Value = 335.79602
$ws.Range("O6").Value = 0.8218100075305903
$ws.Range("P6").Value = 0.8218100075305903
$ws.Range("Q6").Value = 2210.230931895949
$ws.Range("R6").Value = 19892.07838706354
$ws.Range("S6").Value = 0.6726666273951438
$ws.Range("T6").Value = 0.6726666273951438

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 19.74619233333334
$ws.Range("H7").Value = 59.23857700000001
$ws.Range("I7").Value = 0.8185184181638298
$ws.Range("J7").Value = 0.8185184181638298
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 23.563205
$ws.Range("N7").Value = 70.689615
$ws.Range("O7").Value = 0.1730021488506163
$ws.Range("P7").Value = 0.1730021488506163
$ws.Range("Q7").Value = 465.2835779197617
$ws.Range("R7").Value = 4187.552201277856
$ws.Range("S7").Value = 0.1416054452161499
$ws.Range("T7").Value = 0.1416054452161499

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 3.630999
$ws.Range("H8").Value = 10.892997
$ws.Range("I8").Value = 0.1505120332904577
$ws.Range("J8").Value = 0.1505120332904577
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 0.7065936666666666
$ws.Range("N8").Value = 2.119781
$ws.Range("O8").Value = 0.005187843618793344
$ws.Range("P8").Value = 0.005187843618793344
$ws.Range("Q8").Value = 2.565640897073
$ws.Range("R8").Value = 23.090768073657
$ws.Range("S8").Value = 0.0007808328914575122
$ws.Range("T8").Value = 0.0007808328914575124

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 3.630999
$ws.Range("H9").Value = 10.892997
$ws.Range("I9").Value = 0.1505120332904577
$ws.Range("J9").Value = 0.1505120332904577
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 111.9320066666667
$ws.Range("N9").Value = 335.79602
$ws.Range("O9").Value = 0.8218100075305903
$ws.Range("P9").Value = 0.8218100075305903
$ws.Range("Q9").Value = 406.42500427466
$ws.Range("R9").Value = 3657.82503847194
$ws.Range("S9").Value = 0.1236922952118755
$ws.Range("T9").Value = 0.1236922952118755

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 3.630999
$ws.Range("H10").Value = 10.892997
$ws.Range("I10").Value = 0.1505120332904577
$ws.Range("J10").Value = 0.1505120332904577
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 23.563205
$ws.Range("N10").Value = 70.689615
$ws.Range("O10").Value = 0.1730021488506163
$ws.Range("P10").Value = 0.1730021488506163
$ws.Range("Q10").Value = 85.55797379179499
$ws.Range("R10").Value = 770.021764126155
$ws.Range("S10").Value = 0.02603890518712467
$ws.Range("T10").Value = 0.02603890518712468
